$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = '26.089.53'
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,5).NumberFormat = "@"
$ws.Cells.Item(2,5).Value = '  +1.00%  '
$ws.Cells.Item(2,5).Style = "Normal"

# Row 3
$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value = '1.750.07'
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).NumberFormat = "@"
$ws.Cells.Item(3,5).Value = '  +0.50%  '
$ws.Cells.Item(3,5).Style = "Normal"

# Row 4
$ws.Cells.Item(4,5).NumberFormat = "@"
$ws.Cells.Item(4,5).Value = '  -0.16%  '
$ws.Cells.Item(4,5).Style = "Normal"

# Row 5
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = '234.63'
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).NumberFormat = "@"
$ws.Cells.Item(5,5).Value = '  +4.19%  '
$ws.Cells.Item(5,5).Style = "Normal"

# Row 6
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = '0.9990'
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).NumberFormat = "@"
$ws.Cells.Item(6,5).Value = '  -0.05%  '
$ws.Cells.Item(6,5).Style = "Normal"

# Row 7
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = '0.5290'
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).NumberFormat = "@"
$ws.Cells.Item(7,5).Value = '  +2.81%  '
$ws.Cells.Item(7,5).Style = "Normal"

# Row 8
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = '0.2798'
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).NumberFormat = "@"
$ws.Cells.Item(8,5).Value = '  +1.22%  '
$ws.Cells.Item(8,5).Style = "Normal"

# Row 9
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = '0.06185'
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).NumberFormat = "@"
$ws.Cells.Item(9,5).Value = '  +1.56%  '
$ws.Cells.Item(9,5).Style = "Normal"

# Row 10
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = '1.745.93'
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).NumberFormat = "@"
$ws.Cells.Item(10,5).Value = '  +0.35%  '
$ws.Cells.Item(10,5).Style = "Normal"

# Row 11
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = '0.07188'
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).NumberFormat = "@"
$ws.Cells.Item(11,5).Value = '  +2.81%  '
$ws.Cells.Item(11,5).Style = "Normal"

# Row 12
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = '15.41'
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).NumberFormat = "@"
$ws.Cells.Item(12,5).Value = '  +1.41%  '
$ws.Cells.Item(12,5).Style = "Normal"

# Row 13
$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = '0.6451'
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).NumberFormat = "@"
$ws.Cells.Item(13,5).Value = '  +2.02%  '
$ws.Cells.Item(13,5).Style = "Normal"

# Row 14
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = '4.618'
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).NumberFormat = "@"
$ws.Cells.Item(14,5).Value = '  +2.63%  '
$ws.Cells.Item(14,5).Style = "Normal"

# Row 15
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = '78.39'
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).NumberFormat = "@"
$ws.Cells.Item(15,5).Value = '  +2.45%  '
$ws.Cells.Item(15,5).Style = "Normal"

# Row 16
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = '0.9992'
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).NumberFormat = "@"
$ws.Cells.Item(16,5).Value = '  -0.08%  '
$ws.Cells.Item(16,5).Style = "Normal"

# Row 17
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = '0.9989'
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,5).NumberFormat = "@"
$ws.Cells.Item(17,5).Value = '  -0.12%  '
$ws.Cells.Item(17,5).Style = "Normal"

# Row 18
$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = '25.992.24'
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).NumberFormat = "@"
$ws.Cells.Item(18,5).Value = '  +0.56%  '
$ws.Cells.Item(18,5).Style = "Normal"

# Row 19
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = '11.67'
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).NumberFormat = "@"
$ws.Cells.Item(19,5).Value = '  +2.00%  '
$ws.Cells.Item(19,5).Style = "Normal"

# Row 20
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = '0.000006722'
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).NumberFormat = "@"
$ws.Cells.Item(20,5).Value = '  +1.47%  '
$ws.Cells.Item(20,5).Style = "Normal"

# Row 21
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = '1.969.44'
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).NumberFormat = "@"
$ws.Cells.Item(21,5).Value = '  +0.57%  '
$ws.Cells.Item(21,5).Style = "Normal"

# Row 22
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = '4.307'
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).NumberFormat = "@"
$ws.Cells.Item(22,5).Value = '  +5.52%  '
$ws.Cells.Item(22,5).Style = "Normal"

# Row 23
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = '8.767'
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).NumberFormat = "@"
$ws.Cells.Item(23,5).Value = '  +3.34%  '
$ws.Cells.Item(23,5).Style = "Normal"

# Row 24
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = '5.228'
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).NumberFormat = "@"
$ws.Cells.Item(24,5).Value = '  +2.58%  '
$ws.Cells.Item(24,5).Style = "Normal"

# Row 25
$ws.Cells.Item(25,5).NumberFormat = "@"
$ws.Cells.Item(25,5).Value = '  +0.95%  '
$ws.Cells.Item(25,5).Style = "Normal"

# Row 26
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = '1.505'
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).NumberFormat = "@"
$ws.Cells.Item(26,5).Value = '  +0.11%  '
$ws.Cells.Item(26,5).Style = "Normal"

# Row 27
$ws.Cells.Item(27,5).NumberFormat = "@"
$ws.Cells.Item(27,5).Value = '  +2.22%  '
$ws.Cells.Item(27,5).Style = "Normal"

# Row 28
$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = '1.804'
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Cells.Item(28,5).NumberFormat = "@"
$ws.Cells.Item(28,5).Value = '  -0.55%  '
$ws.Cells.Item(28,5).Style = "Normal"

# Row 29
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = '104.58'
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(29,5).NumberFormat = "@"
$ws.Cells.Item(29,5).Value = '  +1.86%  '
$ws.Cells.Item(29,5).Style = "Normal"

# Row 30
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = '0.08286'
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(30,5).NumberFormat = "@"
$ws.Cells.Item(30,5).Value = '  +0.29%  '
$ws.Cells.Item(30,5).Style = "Normal"

# Row 31
$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = '3.796'
$ws.Cells.Item(31,4).Style = "Normal"
$ws.Cells.Item(31,5).NumberFormat = "@"
$ws.Cells.Item(31,5).Value = '  +5.30%  '
$ws.Cells.Item(31,5).Style = "Normal"

# Row 32
$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = '3.671'
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Cells.Item(32,5).NumberFormat = "@"
$ws.Cells.Item(32,5).Value = '  +8.35%  '
$ws.Cells.Item(32,5).Style = "Normal"

# Row 33
$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = '0.04560'
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Cells.Item(33,5).NumberFormat = "@"
$ws.Cells.Item(33,5).Value = '  +3.73%  '
$ws.Cells.Item(33,5).Style = "Normal"

# Row 34
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = '2.643'
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(34,5).NumberFormat = "@"
$ws.Cells.Item(34,5).Value = '  +1.04%  '
$ws.Cells.Item(34,5).Style = "Normal"

# Row 35
$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = '1.004'
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Cells.Item(35,5).NumberFormat = "@"
$ws.Cells.Item(35,5).Value = '  +3.68%  '
$ws.Cells.Item(35,5).Style = "Normal"

# Row 36
$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = '0.6341'
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Cells.Item(36,5).NumberFormat = "@"
$ws.Cells.Item(36,5).Value = '  +6.38%  '
$ws.Cells.Item(36,5).Style = "Normal"

# Row 37
$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = '2.709'
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Cells.Item(37,5).NumberFormat = "@"
$ws.Cells.Item(37,5).Value = '  +2.34%  '
$ws.Cells.Item(37,5).Style = "Normal"

# Row 38
$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = '0.01593'
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(38,5).NumberFormat = "@"
$ws.Cells.Item(38,5).Value = '  +3.01%  '
$ws.Cells.Item(38,5).Style = "Normal"

# Row 39
$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = '1.936'
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).NumberFormat = "@"
$ws.Cells.Item(39,5).Value = '  +0.86%  '
$ws.Cells.Item(39,5).Style = "Normal"

# Row 40
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = '0.9987'
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).NumberFormat = "@"
$ws.Cells.Item(40,5).Value = '  -0.03%  '
$ws.Cells.Item(40,5).Style = "Normal"

# Row 41
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = '99.13'
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).NumberFormat = "@"
$ws.Cells.Item(41,5).Value = '  -1.52%  '
$ws.Cells.Item(41,5).Style = "Normal"

# Row 42
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = '0.3924'
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).NumberFormat = "@"
$ws.Cells.Item(42,5).Value = '  +2.92%  '
$ws.Cells.Item(42,5).Style = "Normal"

# Row 43
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = '0.7449'
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).NumberFormat = "@"
$ws.Cells.Item(43,5).Value = '  +2.33%  '
$ws.Cells.Item(43,5).Style = "Normal"

# Row 44
$ws.Cells.Item(44,5).NumberFormat = "@"
$ws.Cells.Item(44,5).Value = '  +3.22%  '
$ws.Cells.Item(44,5).Style = "Normal"

# Row 45
$ws.Cells.Item(45,5).NumberFormat = "@"
$ws.Cells.Item(45,5).Value = '  +4.50%  '
$ws.Cells.Item(45,5).Style = "Normal"

# Row 46
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = '6.333'
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).NumberFormat = "@"
$ws.Cells.Item(46,5).Value = '  +1.39%  '
$ws.Cells.Item(46,5).Style = "Normal"

# Row 47
$ws.Cells.Item(47,5).NumberFormat = "@"
$ws.Cells.Item(47,5).Value = '  -2.30%  '
$ws.Cells.Item(47,5).Style = "Normal"

# Row 48
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = '54.10'
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(48,5).NumberFormat = "@"
$ws.Cells.Item(48,5).Value = '  +4.00%  '
$ws.Cells.Item(48,5).Style = "Normal"

# Row 49
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = '30.80'
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).NumberFormat = "@"
$ws.Cells.Item(49,5).Value = '  +4.05%  '
$ws.Cells.Item(49,5).Style = "Normal"

# Row 50
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = '7.658'
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).NumberFormat = "@"
$ws.Cells.Item(50,5).Value = '  +2.18%  '
$ws.Cells.Item(50,5).Style = "Normal"

# Row 51
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = '0.3460'
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).NumberFormat = "@"
$ws.Cells.Item(51,5).Value = '  +2.28%  '
$ws.Cells.Item(51,5).Style = "Normal"
